$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table originally listed line1..line6 followed by extr1..extr8 (rows
# 2-15). Two new line entries (line7, line8) are inserted into the
# sequence: what used to be "extr1" (row 8) becomes "line7", and what used
# to be "extr2" (row 9) becomes "line8". Every following "extr" entry
# shifts down by two rows, so two brand-new rows (16 and 17) are appended
# at the bottom to hold "extr7" and "extr8".

# Relabel column B from the bottom up so we never overwrite a value before
# using it.
$ws.Range("B17").Value = "extr8"
$ws.Range("B16").Value = "extr7"
$ws.Range("B15").Value = "extr6"
$ws.Range("B14").Value = "extr5"
$ws.Range("B13").Value = "extr4"
$ws.Range("B12").Value = "extr3"
$ws.Range("B11").Value = "extr2"
$ws.Range("B10").Value = "extr1"
$ws.Range("B9").Value  = "line8"
$ws.Range("B8").Value  = "line7"

# New rows 16 and 17 need the same formatting as the rest of column A
# (bold, bordered, centered style).
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null

$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15

# Update the numeric/boolean data for every affected row (8 through 17).
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $false

$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $false

$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false

$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $false

$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false
